$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: fill an empty table cell with bold text (sz 24 / szCs 24),
# matching the formatting already carried by the cell's empty paragraph mark.
function Fill-Cell($row, $col, $value) {
    $cell = $t.Rows.Item($row).Cells.Item($col)
    $rng = $cell.Range
    $rng.Text = $value
    $rng2 = $cell.Range
    $rng2.Font.Bold = 1
    $rng2.Font.Size = 12
    $rng2.Font.SizeBi = 12
}

# RETENTION section: "Ratio" row -> 0.4
Fill-Cell 24 2 "0.4"

# QUESTION AND ANSWER TASK section:
# "Answer Recall Lenient (ARL)" row -> 0.25
Fill-Cell 44 2 "0.25"
# "Answer Recall Strict (ARS)" row -> 0.1666
Fill-Cell 45 2 "0.1666"
# "Answer Recall Average (ARA)" row -> 0.2083
Fill-Cell 46 2 "0.2083"
